# feat: add edge 0-shot
# Adds two summary rows below the existing per-column AVERAGE row (105):
#   row 106 -> STDEV.S for each metric column (B:F)
#   row 107 -> 95% confidence-interval half-width: stdev/SQRT(103)*1.96
# Also moves the visible selection down to reflect the new bottom of the
# sheet (mirrors the author scrolling down after adding the rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 106: sample standard deviation per column --------------------
$ws.Range("B106").Formula = "=STDEV.S(B2:B104)"
$ws.Range("C106:F106").Formula = "=STDEV.S(C2:C104)"

# --- Row 107: 95% CI half-width off the row-106 stdevs ----------------
$ws.Range("B107").Formula = "=B106/SQRT(103)*1.96"
$ws.Range("C107:F107").Formula = "=C106/SQRT(103)*1.96"

# --- View state: scroll down and select near the new bottom rows ------
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1
$ws.Range("H114").Select()
